# This script re-applies a cyclic rotation of several observation rows in
# the "Artfynd" sheet: the data (all cell values, across columns A..AY) that
# used to live in one row now lives in a different row. Row numbers / row
# level formatting stay put; only the field values move between rows.
#
# Because several of the moves form cycles (e.g. row 7 <- row 9,
# row 8 <- row 7, row 9 <- row 8) we must snapshot every source row BEFORE
# writing anything, otherwise later writes would clobber data that a
# subsequent step still needs to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row number -> source row number (source row's data ends
# up in the destination row).
$rowMap = @{
    7  = 9
    8  = 7
    9  = 8
    18 = 19
    19 = 20
    20 = 18
    21 = 22
    22 = 21
    28 = 29
    29 = 31
    30 = 28
    31 = 30
    32 = 35
    33 = 34
    34 = 36
    35 = 32
    36 = 33
}

# Columns A (1) .. AY (51) cover every populated field in these rows.
$firstCol = 1
$lastCol = 51

# 1) Snapshot every distinct source row's values first.
$snapshot = @{}
foreach ($srcRow in ($rowMap.Values | Sort-Object -Unique)) {
    $rowData = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowData[$c] = $ws.Cells.Item($srcRow, $c).Value()
    }
    $snapshot[$srcRow] = $rowData
}

# Columns holding literal "YYYY-MM-DD" text (Y=25, AA=27). Excel's COM
# layer auto-parses such strings into date serial numbers when assigned
# normally, which would change both the stored type and the value shown
# in the XML (<v>46062</v> instead of the original text). Force those two
# columns to be written as plain text, then restore the default ("Normal")
# cell style so no stray formatting/style index gets left behind.
$dateTextCols = @(25, 27)

# 2) Now write the snapshotted values into their destination rows.
foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($destRow, $c)
        if ($dateTextCols -contains $c) {
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$c]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$c]
        }
    }
}
